$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = 'ככל,ש,העותר,מערער'
$ws.Range("J2").Value = 'positive'

$ws.Range("I3").Value = 'ההכרעה,מינית,אותה,את,כולל'
$ws.Range("J3").Value = 'positive'

$ws.Range("I4").Value = 'הצהירה,ב,ההצהרה,כוח,נמחקת'
$ws.Range("J4").Value = 'positive'

$ws.Range("I5").Value = 'בקשתה,ברם,השופט,הפנים,מטעמים'
$ws.Range("J5").Value = 'positive'

$ws.Range("I6").Value = 'מן,לכך,בו,חזר,צו'
$ws.Range("J6").Value = 'positive'

$ws.Range("I7").Value = 'אותנו,מדינה,והנחדלים,וכך,מ'
$ws.Range("J7").Value = 'positive'

$ws.Range("I8").Value = 'כנספח,אלה,המיסוי,הוא,משפחתו'
$ws.Range("J8").Value = 'negative'

$ws.Range("I9").Value = 'כוח,עוד,כי,השאיר,גבוהים'
$ws.Range("J9").Value = 'positive'

$ws.Range("I10").Value = 'אנו,קדם,בהמלצת,שכר,זה'
$ws.Range("J10").Value = 'positive'

$ws.Range("I11").Value = 'בלא,נתייתרה,העתירה,להוצאות,היא'
$ws.Range("J11").Value = 'positive'

$ws.Range("I12").Value = 'כמבוקש,העתירה,תימחק,ללא,צו'
$ws.Range("J12").Value = 'positive'

$ws.Range("I13").Value = 'הפנים,החליט,העותרים,ראש,ההחלטה'
$ws.Range("J13").Value = 'positive'

$ws.Range("I14").Value = 'השופט,התקנות,הנזכרת,הבחירה,יש'
$ws.Range("J14").Value = 'positive'

$ws.Range("I15").Value = 'למכור,ביקש,ידי,מעבירות,נשק'
$ws.Range("J15").Value = 'positive'

$ws.Range("I16").Value = 'החזר,יועבר,הרשם,עניין,צו'
$ws.Range("J16").Value = 'positive'

$ws.Range("I17").Value = 'כדין,בגין,מ,כחברת,דרישה'
$ws.Range("J17").Value = 'positive'

$ws.Range("I18").Value = 'האמור,בערעור,נמחק,נוכח,מתייתר'
$ws.Range("J18").Value = 'positive'

$ws.Range("I19").Value = 'דוחים,בנפרד,העתירה,נימוקינו,יינתנו'

$ws.Range("I20").Value = 'כבוד,משפחה,המערער,וסעיף,בר'
$ws.Range("J20").Value = 'positive'

$ws.Range("I21").Value = 'נדחית,הופסקו,בעתירה,ואין,לנשים'
$ws.Range("J21").Value = 'positive'

$ws.Range("I22").Value = 'צו,נמחקת,להחזר,אגרה,התייתר'
$ws.Range("J22").Value = 'positive'

$ws.Range("I23").Value = 'צו,העתירה,את,לבקשת,נמחקת'
$ws.Range("J23").Value = 'positive'

$ws.Range("I24").Value = 'עיקר,יפוג,יחדשו,עד,במהלך'
$ws.Range("J24").Value = 'positive'

$ws.Range("I25").Value = 'שנים,מאסר,בפועל,של,מדינת'
$ws.Range("J25").Value = 'positive'

$ws.Range("I26").Value = 'אחר,והמנוח,דקירה,מאסר,בלבו'
$ws.Range("J26").Value = 'positive'

$ws.Range("I27").Value = 'בחנו,התכנון,של,למותר,הדין'
$ws.Range("J27").Value = 'positive'

$ws.Range("I28").Value = 'לבקשה,השבוע,להשתחרר,ביניים,אזרחיים'
$ws.Range("J28").Value = 'positive'

$ws.Range("I29").Value = 'מספר,לרישום,אם,האלף,כינוי'
$ws.Range("J29").Value = 'positive'

$ws.Range("I30").Value = 'בעניינו,השופט,עמדת,בית,ופסקי'
$ws.Range("J30").Value = 'positive'

$ws.Range("I31").Value = 'הגיש,בית,התביעה,בניגוד,היכן'
$ws.Range("J31").Value = 'positive'

$ws.Range("I32").Value = 'שנים,זה,בעלי,בין,משפטי'
$ws.Range("J32").Value = 'positive'

$ws.Range("I33").Value = 'וערעור,להגיש,מיום,קבע,מרזל'
$ws.Range("J33").Value = 'positive'

$ws.Range("I34").Value = 'מסוכן,בא,החזקת,חומרתן,זה'
$ws.Range("J34").Value = 'positive'

$ws.Range("I35").Value = 'צו,המשפט,מבוטלים,המחוזי,בית'
$ws.Range("J35").Value = 'positive'

$ws.Range("I36").Value = 'עצמית,למערער,עליו,למשך,זה'
$ws.Range("J36").Value = 'positive'

$ws.Range("I37").Value = 'לבית,האזורי,בתל,בקשה,הגיש'
$ws.Range("J37").Value = 'positive'

$ws.Range("I38").Value = 'על,עד,לפנינו,מקרה,כלפיו'
$ws.Range("J38").Value = 'positive'

$ws.Range("I39").Value = 'במשרד,שנים,מחש,במשטרת,אין'
$ws.Range("J39").Value = 'positive'

$ws.Range("I40").Value = 'שוחרר,המערער,עונש,שתי,לשנים'
$ws.Range("J40").Value = 'positive'

$ws.Range("I41").Value = 'יחדש,ברשות,ששה,חידוש,הרכב'
$ws.Range("J41").Value = 'positive'

$ws.Range("I42").Value = 'של,הגמל,קשר,טוענים,הכרח'
$ws.Range("J42").Value = 'positive'
